$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.860.65"
$ws.Range("E2").Value = "  +3.36%  "

$ws.Range("D3").Value = "3.980.65"
$ws.Range("E3").Value = "  +1.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.13"
$ws.Range("E5").Value = "  +9.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.52"
$ws.Range("E6").Value = "  +7.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.683"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.750"
$ws.Range("E9").Value = "  +2.95%  "

$ws.Range("E10").Value = "  +1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.77"
$ws.Range("E11").Value = "  +2.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000321"
$ws.Range("E12").Value = "  +2.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.83"
$ws.Range("E13").Value = "  +3.75%  "

$ws.Range("D14").Value = "4.614.21"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("D15").Value = "3.985.28"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("E16").Value = "  +9.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.07"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.38"
$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").Value = "72.502.65"
$ws.Range("E20").Value = "  +2.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.81"
$ws.Range("E21").Value = "  +2.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.68"
$ws.Range("E22").Value = "  +12.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.21"
$ws.Range("E23").Value = "  +0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.43"
$ws.Range("E24").Value = "  -1.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.34"
$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("E26").Value = "  +23.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.08"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.68"
$ws.Range("E28").Value = "  +2.42%  "

$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.48"
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  +6.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.62"
$ws.Range("E32").Value = "  +2.51%  "

$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "679.66"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.54"
$ws.Range("E35").Value = "  +2.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "68.89"
$ws.Range("E36").Value = "  +6.24%  "

$ws.Range("D37").Value = "0.0₃0880"
$ws.Range("E37").Value = "  +8.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.437"
$ws.Range("E38").Value = "  +2.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("E42").Value = "  -2.08%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0487"
$ws.Range("E44").Value = "  +2.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.79"
$ws.Range("E45").Value = "  +12.70%  "

$ws.Range("E46").Value = "  +1.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.66"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.39"
$ws.Range("E48").Value = "  +2.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.04"
$ws.Range("E49").Value = "  +3.21%  "

$ws.Range("E50").Value = "  +6.13%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.146"
$ws.Range("E39").Value = "  -0.39%  "

# Row 40
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.35"
$ws.Range("E40").Value = "  +1.70%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  +8.72%  "

